$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.998.39'
$ws.Range("E2").Value = '  -1.96%  '
$ws.Range("D3").Value = '1.883.11'
$ws.Range("E3").Value = '  -1.53%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.77%  '
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5000'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.38%  '
$ws.Range("E8").Value = '  -2.72%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09159'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.47%  '
$ws.Range("E10").Value = '  -2.76%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.66'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.86%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.328'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.07%  '
$ws.Range("E13").Value = '  -2.46%  '
$ws.Range("D14").Value = '1.881.30'
$ws.Range("E14").Value = '  -1.69%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.279'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001104'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.74%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.40'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.58%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06640'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.97'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.32%  '
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.181'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.01%  '
$ws.Range("D23").Value = '28.047.01'
$ws.Range("E23").Value = '  -2.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.38'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.301'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.60%  '
$ws.Range("D26").Value = '2.098.30'
$ws.Range("E26").Value = '  -1.62%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.547'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.80%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.76'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.38%  '
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '157.55'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.55%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.64'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.74%  '
$ws.Range("E31").Value = '  -2.51%  '
$ws.Range("E32").Value = '  -3.94%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.591'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.582'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.338'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.70%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06581'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02404'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2190'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.59%  '
$ws.Range("E39").Value = '  +8.07%  '
$ws.Range("E40").Value = '  -5.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6405'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.54'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.935'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9997'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.28'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.59%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6035'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.292'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.670'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.990'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.213'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.70'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.57%  '
